# The page-number reminder paragraph currently reads "118".
# Update it to the new page number "120".
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("118", $true, $false, $false, $false, $false, $true, 1, $false, "120", 2)
